$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "56.659.42"
Set-TextValue "E2" "  +0.37%  "
Set-TextValue "D3" "3.022.93"
Set-TextValue "E3" "  +2.69%  "
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "510.79"
Set-TextValue "E5" "  +3.57%  "
Set-TextValue "D6" "139.94"
Set-TextValue "E6" "  +4.84%  "
Set-TextValue "E7" "  -0.06%  "
Set-TextValue "D8" "0.431"
Set-TextValue "E8" "  +2.25%  "
Set-TextValue "D9" "7.12"
Set-TextValue "E9" "  +0.46%  "
Set-TextValue "E10" "  +2.81%  "
Set-TextValue "E11" "  +5.88%  "
Set-TextValue "D12" "3.545.68"
Set-TextValue "E12" "  +2.59%  "
Set-TextValue "E13" "  +1.02%  "
Set-TextValue "D14" "25.34"
Set-TextValue "E14" "  -1.67%  "
Set-TextValue "E15" "  +4.33%  "
Set-TextValue "D16" "56.656.37"
Set-TextValue "E16" "  +0.22%  "
Set-TextValue "D17" "3.028.31"
Set-TextValue "E17" "  +3.02%  "
Set-TextValue "D18" "5.91"
Set-TextValue "E18" "  -0.29%  "
Set-TextValue "D19" "13.11"
Set-TextValue "E19" "  +6.05%  "
Set-TextValue "E20" "  +4.35%  "
Set-TextValue "D21" "334.21"
Set-TextValue "E21" "  +6.30%  "
Set-TextValue "E22" "  +0.02%  "
Set-TextValue "D23" "0.500"
Set-TextValue "E23" "  +4.07%  "
Set-TextValue "D24" "64.66"
Set-TextValue "E24" "  +3.63%  "
Set-TextValue "D25" "3.156.40"
Set-TextValue "E25" "  +2.70%  "
Set-TextValue "E26" "  +3.65%  "
Set-TextValue "E27" "  -0.18%  "
Set-TextValue "D28" "0.0₃0926"
Set-TextValue "E28" "  +9.11%  "
Set-TextValue "D29" "6.38"
Set-TextValue "E29" "  -0.12%  "
Set-TextValue "E30" "  -2.26%  "
Set-TextValue "E31" "  +3.38%  "
Set-TextValue "D32" "20.42"
Set-TextValue "E32" "  +3.11%  "
Set-TextValue "E33" "  +3.35%  "
Set-TextValue "D34" "152.86"
Set-TextValue "E34" "  +0.95%  "
Set-TextValue "D35" "4.48"
Set-TextValue "E35" "  +1.38%  "
Set-TextValue "E36" "  +15.31%  "
Set-TextValue "E37" "  +3.00%  "
Set-TextValue "E38" "  +2.56%  "
Set-TextValue "E39" "  +2.20%  "
Set-TextValue "D40" "3.062.81"
Set-TextValue "E40" "  +2.96%  "
Set-TextValue "D41" "36.47"
Set-TextValue "E42" "  -0.06%  "
Set-TextValue "D43" "3.80"
Set-TextValue "E43" "  +4.29%  "
Set-TextValue "E44" "  +3.87%  "
Set-TextValue "D45" "2.214.02"
Set-TextValue "E45" "  +3.94%  "
Set-TextValue "E46" "  +1.16%  "
Set-TextValue "E47" "  +6.39%  "
Set-TextValue "D48" "0.931"
Set-TextValue "E48" "  +2.05%  "
Set-TextValue "E49" "  +5.69%  "
Set-TextValue "D50" "5.82"
Set-TextValue "E50" "  +0.15%  "
Set-TextValue "D51" "0.0856"
Set-TextValue "E51" "  +1.67%  "
